$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Ajout de noms et prenoms fictifs (colonnes A = Nom, C = Prenom) pour les lignes 4 a 41
$ws.Range("A4").Value = "THOMAS"
$ws.Range("C4").Value = "Michel"
$ws.Range("A5").Value = "PETIT"
$ws.Range("C5").Value = "Pierre"
$ws.Range("A6").Value = "ROBERT"
$ws.Range("C6").Value = "Philippe"
$ws.Range("A7").Value = "RICHARD"
$ws.Range("C7").Value = "Alain"
$ws.Range("A8").Value = "DURAND"
$ws.Range("C8").Value = "Nathalie"
$ws.Range("A9").Value = "DUBOIS"
$ws.Range("C9").Value = "Dominique"
$ws.Range("A10").Value = "MOREAU"
$ws.Range("C10").Value = "Nicolas"
$ws.Range("A11").Value = "LAURENT"
$ws.Range("C11").Value = "Isabelle"
$ws.Range("A12").Value = "SIMON"
$ws.Range("C12").Value = "Catherine"
$ws.Range("A13").Value = "MICHEL"
$ws.Range("C13").Value = "Monique"
$ws.Range("A14").Value = "LEFEBVRE"
$ws.Range("C14").Value = "Sylvie"
$ws.Range("A15").Value = "LEROY"
$ws.Range("C15").Value = "Patrick"
$ws.Range("A16").Value = "ROUX"
$ws.Range("C16").Value = "Bernard"
$ws.Range("A17").Value = "DAVID"
$ws.Range("C17").Value = "Françoise"
$ws.Range("A18").Value = "BERTRAND"
$ws.Range("C18").Value = "Claude"
$ws.Range("A19").Value = "MOREL"
$ws.Range("C19").Value = "Daniel"
$ws.Range("A20").Value = "FOURNIER"
$ws.Range("C20").Value = "Christophe"
$ws.Range("A21").Value = "GIRARD"
$ws.Range("C21").Value = "Christian"
$ws.Range("A22").Value = "BONNET"
$ws.Range("C22").Value = "André"
$ws.Range("A23").Value = "DUPONT"
$ws.Range("C23").Value = "Jacques"
$ws.Range("A24").Value = "LAMBERT"
$ws.Range("C24").Value = "Martine"
$ws.Range("A25").Value = "FONTAINE"
$ws.Range("C25").Value = "Gérard"
$ws.Range("A26").Value = "ROUSSEAU"
$ws.Range("C26").Value = "Jacqueline"
$ws.Range("A27").Value = "VINCENT"
$ws.Range("C27").Value = "Frédéric"
$ws.Range("A28").Value = "MULLER"
$ws.Range("C28").Value = "Éric"
$ws.Range("A29").Value = "LEFEVRE"
$ws.Range("C29").Value = "Laurent"
$ws.Range("A30").Value = "FAURE"
$ws.Range("C30").Value = "Julien"
$ws.Range("A31").Value = "ANDRE"
$ws.Range("C31").Value = "David"
$ws.Range("A32").Value = "MERCIER"
$ws.Range("C32").Value = "Stéphane"
$ws.Range("A33").Value = "BLANC"
$ws.Range("C33").Value = "Sébastien"
$ws.Range("A34").Value = "GUERIN"
$ws.Range("C34").Value = "Anne"
$ws.Range("A35").Value = "BOYER"
$ws.Range("C35").Value = "Pascal"
$ws.Range("A36").Value = "GARNIER"
$ws.Range("C36").Value = "Christine"
$ws.Range("A37").Value = "CHEVALIER"
$ws.Range("C37").Value = "Nicole"
$ws.Range("A38").Value = "FRANCOIS"
$ws.Range("C38").Value = "Thierry"
$ws.Range("A39").Value = "LEGRAND"
$ws.Range("C39").Value = "Olivier"
$ws.Range("A40").Value = "GAUTHIER"
$ws.Range("C40").Value = "Thomas"
$ws.Range("A41").Value = "GARCIA"
$ws.Range("C41").Value = "Alexandre"

# Harmoniser le format/style de la colonne V (statut) sur les lignes 2 a 41 : passage au format "000"
for ($r = 2; $r -le 41; $r++) {
    $ws.Range("V$r").NumberFormat = "000"
}

# Les lignes 29, 31 et 38 portaient une police rouge distincte : on recale leur mise en forme
# sur celle (deja normalisee) de V2 pour retrouver exactement le meme style.
$ws.Range("V2").Copy()
$ws.Range("V29").PasteSpecial(-4122)
$ws.Range("V31").PasteSpecial(-4122)
$ws.Range("V38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Largeurs de colonnes A/B/C ajustees (saisie des nouvelles donnees) et cellule active finale
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 14.1
$ws.Columns.Item(3).ColumnWidth = 13.6
$ws.Range("B18").Select()

